# Add new columns I (I0) and J (IF) to the sheet, mirroring the
# existing header style used by columns B..H, and fill in the data
# for rows 2..37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font + border) from the existing
# header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (rows 2..37) --------------------------------------------
# Each tuple is: row, I-value, J-value
# NOTE: each inner @(...) is prefixed with the unary comma operator so
# the outer @(...) keeps it as a nested array element instead of
# flattening everything into one long list.
$data = @(
    ,@(2, 5, 5)
    ,@(3, 7, 8)
    ,@(4, 1, 5)
    ,@(5, 1, 6)
    ,@(6, 1, 5)
    ,@(7, 1, 6)
    ,@(8, 1, 9)
    ,@(9, 1, 6)
    ,@(10, 1, 4)
    ,@(11, 1, 7)
    ,@(12, 1, 5)
    ,@(13, 5, 6)
    ,@(14, 1, 4)
    ,@(15, 1, 5)
    ,@(16, 1, 7)
    ,@(17, 1, 5)
    ,@(18, 1, 6)
    ,@(19, 1, 6)
    ,@(20, 1, 6)
    ,@(21, 4, 8)
    ,@(22, 1, 5)
    ,@(23, 1, 5)
    ,@(24, 1, 3)
    ,@(25, 1, 6)
    ,@(26, 1, 6)
    ,@(27, 1, 6)
    ,@(28, 1, 5)
    ,@(29, 1, 5)
    ,@(30, 1, 6)
    ,@(31, 1, 5)
    ,@(32, 1, 6)
    ,@(33, 1, 4)
    ,@(34, 1, 4)
    ,@(35, 1, 3)
    ,@(36, 1, 2)
    ,@(37, 1, 1)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
